$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9476458981920928
$ws.Range("C2").Value = 0.1719206942726146
$ws.Range("D2").Value = 0.1443690170295326
$ws.Range("E2").Value = 0.09112395411323604
$ws.Range("F2").Value = 3.037772572775197
$ws.Range("M2").Value = 0.3725926983130634
$ws.Range("B3").Value = 0.8641610323120972
$ws.Range("C3").Value = 0.1504687380920871
$ws.Range("D3").Value = 0.1336508942510051
$ws.Range("E3").Value = 0.08435741998354018
$ws.Range("F3").Value = 2.820701146081149
$ws.Range("M3").Value = 0.3392264900062685
$ws.Range("B4").Value = 0.8137861859623285
$ws.Range("C4").Value = 0.1374226980942126
$ws.Range("D4").Value = 0.1270718285175718
$ws.Range("E4").Value = 0.08024515227905482
$ws.Range("F4").Value = 2.688233000857053
$ws.Range("M4").Value = 0.3190326487632902
$ws.Range("B5").Value = 0.7934772783991946
$ws.Range("C5").Value = 0.132136802423247
$ws.Range("D5").Value = 0.1243906508101418
$ws.Range("E5").Value = 0.0785796964358525
$ws.Range("F5").Value = 2.634446044827286
$ws.Range("M5").Value = 0.3108757524015218
$ws.Range("B6").Value = 0.7901181458481688
$ws.Range("C6").Value = 0.1312608921788012
$ws.Range("D6").Value = 0.1239454154905673
$ws.Range("E6").Value = 0.07830376266257844
$ws.Range("F6").Value = 2.625526243975685
$ws.Range("M6").Value = 0.3095256281440086
$ws.Range("B7").Value = 0.8135114095486244
$ws.Range("C7").Value = 0.1373512887144557
$ws.Range("D7").Value = 0.1270356706801294
$ws.Range("E7").Value = 0.08022264999094375
$ws.Range("F7").Value = 2.687506835130364
$ws.Range("M7").Value = 0.318922351340845
$ws.Range("B8").Value = 0.9186745346849534
$ws.Range("C8").Value = 0.1644973286932725
$ws.Range("D8").Value = 0.1406726273246051
$ws.Range("E8").Value = 0.08878184398057698
$ws.Range("F8").Value = 2.962751197177795
$ws.Range("M8").Value = 0.3610262927193446
$ws.Range("B9").Value = 1.132087899232545
$ws.Range("C9").Value = 0.218779793981895
$ws.Range("D9").Value = 0.1674595635140577
$ws.Range("E9").Value = 0.1059191192272593
$ws.Range("F9").Value = 3.509445478715321
$ws.Range("M9").Value = 0.4459890776067965
$ws.Range("B10").Value = 1.293515007262158
$ws.Range("C10").Value = 0.2593791426511416
$ws.Range("D10").Value = 0.187210288750066
$ws.Range("E10").Value = 0.118748732103505
$ws.Range("F10").Value = 3.916036404187423
$ws.Range("M10").Value = 0.50997893781269
$ws.Range("B11").Value = 1.368014845915297
$ws.Range("C11").Value = 0.278022386117442
$ws.Range("D11").Value = 0.1962201569249089
$ws.Range("E11").Value = 0.1246424250082967
$ws.Range("F11").Value = 4.102233106722736
$ws.Range("M11").Value = 0.5394539449726494
$ws.Range("B12").Value = 1.396383658376237
$ws.Range("C12").Value = 0.2851085723090137
$ws.Range("D12").Value = 0.1996363256629934
$ws.Range("E12").Value = 0.1268828873086179
$ws.Range("F12").Value = 4.17293092655774
$ws.Range("M12").Value = 0.5506698174898474
$ws.Range("B13").Value = 1.390266863196246
$ws.Range("C13").Value = 0.2835812401041835
$ws.Range("D13").Value = 0.1989003887631213
$ws.Range("E13").Value = 0.1263999727388665
$ws.Range("F13").Value = 4.157696289548142
$ws.Range("M13").Value = 0.5482518343006717
$ws.Range("B14").Value = 1.370345588116095
$ws.Range("C14").Value = 0.2786048342753702
$ws.Range("D14").Value = 0.1965011163645443
$ws.Range("E14").Value = 0.1248265736015455
$ws.Range("F14").Value = 4.108045601822823
$ws.Range("M14").Value = 0.5403755826417722
$ws.Range("B15").Value = 1.358163845118895
$ws.Range("C15").Value = 0.2755601190013124
$ws.Range("D15").Value = 0.195032077474508
$ws.Range("E15").Value = 0.1238639591919792
$ws.Range("F15").Value = 4.077658098306017
$ws.Range("M15").Value = 0.5355582784262225
$ws.Range("B16").Value = 1.288668060711927
$ws.Range("C16").Value = 0.2581643839289143
$ws.Range("D16").Value = 0.1866220393951892
$ws.Range("E16").Value = 0.1183647543156212
$ws.Range("F16").Value = 3.903893859664578
$ws.Range("M16").Value = 0.5080601833441705
$ws.Range("B17").Value = 1.246310498272294
$ws.Range("C17").Value = 0.2475382572559397
$ws.Range("D17").Value = 0.1814696997097371
$ws.Range("E17").Value = 0.1150061561220852
$ws.Range("F17").Value = 3.797619838831224
$ws.Range("M17").Value = 0.4912858816138339
$ws.Range("B18").Value = 1.222047619935609
$ws.Range("C18").Value = 0.2414427244915203
$ws.Range("D18").Value = 0.1785085358092431
$ws.Range("E18").Value = 0.1130797618720081
$ws.Range("F18").Value = 3.736609521624132
$ws.Range("M18").Value = 0.4816720616151429
$ws.Range("B19").Value = 1.213849691232724
$ws.Range("C19").Value = 0.2393816476850361
$ws.Range("D19").Value = 0.1775063145929181
$ws.Range("E19").Value = 0.1124284306039556
$ws.Range("F19").Value = 3.71597200302142
$ws.Range("M19").Value = 0.4784228286096948
$ws.Range("B20").Value = 1.250809146951724
$ws.Range("C20").Value = 0.2486677266165884
$ws.Range("D20").Value = 0.1820179292715522
$ws.Range("E20").Value = 0.1153631247576996
$ws.Range("F20").Value = 3.808920834708147
$ws.Range("M20").Value = 0.4930679707592702
$ws.Range("B21").Value = 1.376192646304048
$ws.Range("C21").Value = 0.2800657982993187
$ws.Range("D21").Value = 0.197205717692583
$ws.Range("E21").Value = 0.1252884812070789
$ws.Range("F21").Value = 4.122623986738404
$ws.Range("M21").Value = 0.5426875401807365
$ws.Range("B22").Value = 1.459057083136486
$ws.Range("C22").Value = 0.3007407355483736
$ws.Range("D22").Value = 0.2071573736814116
$ws.Range("E22").Value = 0.1318258725213468
$ws.Range("F22").Value = 4.328755178232143
$ws.Range("M22").Value = 0.5754342946179776
$ws.Range("B23").Value = 1.414745303010591
$ws.Range("C23").Value = 0.2896915535639266
$ws.Range("D23").Value = 0.201843412269767
$ws.Range("E23").Value = 0.128331983614693
$ws.Range("F23").Value = 4.218633886059536
$ws.Range("M23").Value = 0.5579270984154476
$ws.Range("B24").Value = 1.248775030824049
$ws.Range("C24").Value = 0.2481570513533882
$ws.Range("D24").Value = 0.1817700716723323
$ws.Range("E24").Value = 0.1152017252236064
$ws.Range("F24").Value = 3.803811381126224
$ws.Range("M24").Value = 0.4922621951016168
$ws.Range("B25").Value = 1.073558865653581
$ws.Range("C25").Value = 0.2039746313485296
$ws.Range("D25").Value = 0.1602037670672161
$ws.Range("E25").Value = 0.1012428306904525
$ws.Range("F25").Value = 3.360735064642626
$ws.Range("M25").Value = 0.422736911962744

Write-Output "Applied 144 cell updates"
